$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44330, 0, 3, 32.31365790607497),
    @(44331, 2, 4, 43.08487720809995),
    @(44332, 0, 4, 43.08487720809995),
    @(44333, 1, 4, 43.08487720809995),
    @(44334, 1, 4, 43.08487720809995),
    @(44335, 0, 4, 43.08487720809995),
    @(44336, 0, 4, 43.08487720809995),
    @(44337, 2, 6, 64.62731581214993),
    @(44338, 2, 6, 64.62731581214993),
    @(44339, 4, 10, 107.7121930202499),
    @(44340, 1, 10, 107.7121930202499),
    @(44341, 0, 9, 96.9409737182249),
    @(44342, 0, 9, 96.9409737182249),
    @(44343, 1, 10, 107.7121930202499)
)

$startRow = 256
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$endRow = $startRow + $data.Count - 1

# Carry the date-column formatting (style) from the last existing row down to the new rows
$ws.Range("A255").Copy()
$ws.Range("A256:A$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
